$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.284.24'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").Value = '1.887.83'
$ws.Range("E3").Value = '  -1.20%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4683'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2860'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06609'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.12'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07791'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '98.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.87%  '
$ws.Range("D13").Value = '1.895.32'
$ws.Range("E13").Value = '  -0.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.129'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6778'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '285.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +12.34%  '
$ws.Range("D17").Value = '30.289.31'
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").Value = '2.143.09'
$ws.Range("E20").Value = '  -0.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.392'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007310'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.13%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.206'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.445'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.27'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09738'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.455'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.487'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.74%  '
$ws.Range("E33").Value = '  -1.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04706'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7130'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.099'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01884'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.706'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.532'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.69'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.82%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8731'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.41%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.978'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '104.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4215'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '987.80'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.276'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.222'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1163'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '34.18'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.63%  '
